# Apply revision data update to S17 Table
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update values in column G (revision data)
$ws.Range("G2").Value = 165
$ws.Range("G4").Value = 145
$ws.Range("G19").Value = 53
$ws.Range("G20").Value = 21
$ws.Range("G21").Value = 5
$ws.Range("G23").Value = 32
$ws.Range("G24").Value = 20

# Update the active selection to match the new state (G2:G25, active cell G2)
$ws.Range("G2:G25").Select()
